$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC!row40
$ws_ALC.Range("H40").Value = 1053.1562
$ws_ALC.Range("I40").Value = 950
$ws_ALC.Range("J40").Value = 1076.9615
$ws_ALC.Range("K40").Value = 950
$ws_ALC.Range("L40").Value = 1076.9615
$ws_ALC.Range("M40").Value = -775
$ws_ALC.Range("N40").Value = -1426.9615

# ALC!row132
$ws_ALC.Range("H132").Value = 1443806.2
$ws_ALC.Range("I132").Value = 2731.625
$ws_ALC.Range("J132").Value = 24501000
$ws_ALC.Range("K132").Value = 8194.875
$ws_ALC.Range("L132").Value = 73503000
$ws_ALC.Range("M132").Value = -5664.875
$ws_ALC.Range("N132").Value = -73508060

# ALC!row141
$ws_ALC.Range("H141").Value = 3065
$ws_ALC.Range("I141").Value = 2095
$ws_ALC.Range("J141").Value = 3550
$ws_ALC.Range("K141").Value = 6285
$ws_ALC.Range("L141").Value = 10650
$ws_ALC.Range("M141").Value = -1105
$ws_ALC.Range("N141").Value = -21010

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM!row32
$ws_ARM.Range("H32").Value = 5072322
$ws_ARM.Range("I32").Value = 6311383
$ws_ARM.Range("J32").Value = 20765.77
$ws_ARM.Range("K32").Value = 6311383
$ws_ARM.Range("L32").Value = 20765.77
$ws_ARM.Range("M32").Value = -6311096
$ws_ARM.Range("N32").Value = -21339.77

# ARM!row45
$ws_ARM.Range("H45").Value = 3162.28
$ws_ARM.Range("I45").Value = 3285.1
$ws_ARM.Range("K45").Value = 3285.1
$ws_ARM.Range("M45").Value = -2908.1

# ARM!row74
$ws_ARM.Range("H74").Value = 7412857.5
$ws_ARM.Range("I74").Value = 12550761
$ws_ARM.Range("K74").Value = 12550761
$ws_ARM.Range("M74").Value = -12549887

# ARM!row77
$ws_ARM.Range("H77").Value = 7412857.5
$ws_ARM.Range("I77").Value = 12550761
$ws_ARM.Range("K77").Value = 62753805
$ws_ARM.Range("M77").Value = -62749437

# ARM!row97
$ws_ARM.Range("H97").Value = 3127157
$ws_ARM.Range("I97").Value = 4168342
$ws_ARM.Range("J97").Value = 3602
$ws_ARM.Range("K97").Value = 4168342
$ws_ARM.Range("L97").Value = 3602
$ws_ARM.Range("M97").Value = -4167846
$ws_ARM.Range("N97").Value = -4594

# ARM!row122
$ws_ARM.Range("H122").Value = 3970140.8
$ws_ARM.Range("I122").Value = 1996.238
$ws_ARM.Range("J122").Value = 15874574
$ws_ARM.Range("K122").Value = 5988.714
$ws_ARM.Range("L122").Value = 47623722
$ws_ARM.Range("M122").Value = -3538.714
$ws_ARM.Range("N122").Value = -47628622

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM!row57
$ws_BSM.Range("H57").Value = 0
$ws_BSM.Range("J57").Value = 0
$ws_BSM.Range("L57").Value = 0
$ws_BSM.Range("N57").ClearContents()

# BSM!row136
$ws_BSM.Range("H136").Value = 0
$ws_BSM.Range("J136").Value = 0
$ws_BSM.Range("L136").Value = 0
$ws_BSM.Range("N136").ClearContents()

# BSM!row140
$ws_BSM.Range("H140").Value = 56555.79
$ws_BSM.Range("J140").Value = 56555.79
$ws_BSM.Range("L140").Value = 56555.79
$ws_BSM.Range("N140").Value = -66915.79000000001

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP!row132
$ws_CRP.Range("H132").Value = 41959.64
$ws_CRP.Range("I132").Value = 1599.65
$ws_CRP.Range("K132").Value = 4798.950000000001
$ws_CRP.Range("M132").Value = -2268.950000000001

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL!row39
$ws_CUL.Range("H39").Value = 2877
$ws_CUL.Range("J39").Value = 2877
$ws_CUL.Range("L39").Value = 8631
$ws_CUL.Range("N39").Value = -9219

# CUL!row40
$ws_CUL.Range("H40").Value = 37.6
$ws_CUL.Range("I40").Value = 37.6
$ws_CUL.Range("J40").Value = 0
$ws_CUL.Range("K40").Value = 150.4
$ws_CUL.Range("L40").Value = 0
$ws_CUL.Range("M40").Value = -81.40000000000001
$ws_CUL.Range("N40").ClearContents()

# CUL!row48
$ws_CUL.Range("H48").Value = 5666.6665
$ws_CUL.Range("J48").Value = 5666.6665
$ws_CUL.Range("L48").Value = 16999.9995
$ws_CUL.Range("N48").Value = -17499.9995

# CUL!row58
$ws_CUL.Range("H58").Value = 1457.1428
$ws_CUL.Range("J58").Value = 1457.1428
$ws_CUL.Range("L58").Value = 4371.428400000001
$ws_CUL.Range("N58").Value = -4627.428400000001

# CUL!row68
$ws_CUL.Range("H68").Value = 808.4796
$ws_CUL.Range("I68").Value = 557.04083
$ws_CUL.Range("J68").Value = 1059.9183
$ws_CUL.Range("K68").Value = 1671.12249
$ws_CUL.Range("L68").Value = 3179.7549
$ws_CUL.Range("M68").Value = -860.1224900000002
$ws_CUL.Range("N68").Value = -4801.7549

# CUL!row70
$ws_CUL.Range("H70").Value = 2769.2307
$ws_CUL.Range("I70").Value = 1550
$ws_CUL.Range("J70").Value = 3311.111
$ws_CUL.Range("K70").Value = 4650
$ws_CUL.Range("L70").Value = 9933.332999999999
$ws_CUL.Range("M70").Value = -4335
$ws_CUL.Range("N70").Value = -10563.333

# CUL!row71
$ws_CUL.Range("H71").Value = 808.4796
$ws_CUL.Range("I71").Value = 557.04083
$ws_CUL.Range("J71").Value = 1059.9183
$ws_CUL.Range("K71").Value = 5013.36747
$ws_CUL.Range("L71").Value = 9539.2647
$ws_CUL.Range("M71").Value = -957.3674700000001
$ws_CUL.Range("N71").Value = -17651.2647

# CUL!row73
$ws_CUL.Range("H73").Value = 2769.2307
$ws_CUL.Range("I73").Value = 1550
$ws_CUL.Range("J73").Value = 3311.111
$ws_CUL.Range("K73").Value = 4650
$ws_CUL.Range("L73").Value = 9933.332999999999
$ws_CUL.Range("M73").Value = -3558
$ws_CUL.Range("N73").Value = -12117.333

# CUL!row76
$ws_CUL.Range("H76").Value = 3377.7778
$ws_CUL.Range("I76").Value = 3000
$ws_CUL.Range("J76").Value = 3425
$ws_CUL.Range("K76").Value = 9000
$ws_CUL.Range("L76").Value = 10275
$ws_CUL.Range("N76").Value = -11041
$ws_CUL.Range("M76").Value = -8617

# CUL!row79
$ws_CUL.Range("H79").Value = 3377.7778
$ws_CUL.Range("I79").Value = 3000
$ws_CUL.Range("J79").Value = 3425
$ws_CUL.Range("K79").Value = 9000
$ws_CUL.Range("L79").Value = 10275
$ws_CUL.Range("N79").Value = -12927
$ws_CUL.Range("M79").Value = -7674

# CUL!row106
$ws_CUL.Range("H106").Value = 3600
$ws_CUL.Range("J106").Value = 3600
$ws_CUL.Range("L106").Value = 10800
$ws_CUL.Range("N106").Value = -12692

# CUL!row109
$ws_CUL.Range("H109").Value = 2153.3845
$ws_CUL.Range("I109").Value = 984
$ws_CUL.Range("J109").Value = 3155.7144
$ws_CUL.Range("K109").Value = 2952
$ws_CUL.Range("L109").Value = 9467.143199999999
$ws_CUL.Range("M109").Value = -1912
$ws_CUL.Range("N109").Value = -11547.1432

# CUL!row129
$ws_CUL.Range("H129").Value = 2138502.8
$ws_CUL.Range("I129").Value = 812.2273
$ws_CUL.Range("J129").Value = 4904925.5
$ws_CUL.Range("K129").Value = 2436.6819
$ws_CUL.Range("L129").Value = 14714776.5
$ws_CUL.Range("M129").Value = 2563.3181
$ws_CUL.Range("N129").Value = -14724776.5

# CUL!row131
$ws_CUL.Range("H131").Value = 819.7879
$ws_CUL.Range("I131").Value = 443.75
$ws_CUL.Range("J131").Value = 940.12
$ws_CUL.Range("K131").Value = 1331.25
$ws_CUL.Range("L131").Value = 2820.36
$ws_CUL.Range("M131").Value = 3708.75
$ws_CUL.Range("N131").Value = -12900.36

# CUL!row132
$ws_CUL.Range("H132").Value = 1548.375
$ws_CUL.Range("J132").Value = 1328.05
$ws_CUL.Range("L132").Value = 11952.45
$ws_CUL.Range("N132").Value = -17012.45

# CUL!row133
$ws_CUL.Range("H133").Value = 4295.3335
$ws_CUL.Range("J133").Value = 6799.3335
$ws_CUL.Range("L133").Value = 20398.0005
$ws_CUL.Range("N133").Value = -30518.0005

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW!row55
$ws_LTW.Range("H55").Value = 357
$ws_LTW.Range("I55").Value = 271.83334
$ws_LTW.Range("J55").Value = 442.16666
$ws_LTW.Range("K55").Value = 271.83334
$ws_LTW.Range("L55").Value = 442.16666
$ws_LTW.Range("M55").Value = -98.83334000000002
$ws_LTW.Range("N55").Value = -788.16666

# LTW!row93
$ws_LTW.Range("H93").Value = 1921.4166
$ws_LTW.Range("I93").Value = 1816.6
$ws_LTW.Range("J93").Value = 1996.2858
$ws_LTW.Range("K93").Value = 1816.6
$ws_LTW.Range("L93").Value = 1996.2858
$ws_LTW.Range("M93").Value = -568.5999999999999
$ws_LTW.Range("N93").Value = -4492.2858

# LTW!row123
$ws_LTW.Range("H123").Value = 27843.545
$ws_LTW.Range("J123").Value = 27843.545
$ws_LTW.Range("L123").Value = 27843.545
$ws_LTW.Range("N123").Value = -37643.545

# LTW!row137
$ws_LTW.Range("H137").Value = 24999.883
$ws_LTW.Range("J137").Value = 24999.883
$ws_LTW.Range("L137").Value = 24999.883
$ws_LTW.Range("N137").Value = -35199.883
